# Weekly update: a new daily record is inserted at row 116 for
# "Hortaliza, Femacal de La Calera - Haba"; all the existing records
# that used to occupy rows 116-146 shift down by one row (to 117-147).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 116, pushing every
# following row (116..146) down to (117..147).
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(116, 1).Value  = 3
$ws.Cells.Item(116, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(116, 3).Value  = "Coquimbo"
$ws.Cells.Item(116, 4).Value  = 44736
$ws.Cells.Item(116, 5).Value  = 5
$ws.Cells.Item(116, 6).Value  = 100112026
$ws.Cells.Item(116, 7).Value  = "Haba"
$ws.Cells.Item(116, 8).Value  = "Sin especificar"
$ws.Cells.Item(116, 9).Value  = "Primera"
$ws.Cells.Item(116, 10).Value = 76
$ws.Cells.Item(116, 11).Value = 21000
$ws.Cells.Item(116, 12).Value = 22000
$ws.Cells.Item(116, 13).Value = 21500
$ws.Cells.Item(116, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(116, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(116, 16).Value = 860
$ws.Cells.Item(116, 17).Value = 25
$ws.Cells.Item(116, 18).Value = "Hortaliza"
